$wb = $excel.ActiveWorkbook

# Rename Sheet1 to AddCustomerTest
$ws = $wb.Worksheets.Item(1)
$ws.Name = "AddCustomerTest"

# Fill in data in the order needed so shared-string indices line up
# with the target workbook (column-major, with a couple of quirks).
$ws.Range("A1").Value = "firstname"
$ws.Range("B1").Value = "lastname"
$ws.Range("C1").Value = "postcode"

$ws.Range("A2").Value = "Tai"
$ws.Range("B2").Value = "Le"

$ws.Range("A3").Value = "Giang"
$ws.Range("B3").Value = "Nguyen"
$ws.Range("C3").Value = "Abcd20"

$ws.Range("C2").Value = "Test1000"

$ws.Range("D1").Value = "alerttext"
$ws.Range("D2").Value = "Customer added successfully"
$ws.Range("D3").Value = "Customer added successfully"

# Update selection to D3 to match the final saved view state
$ws.Range("D3").Select()
